$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 9 of data - a new training run entry
$ws.Range("A9").Value = "regular"
$ws.Range("B9").Value = "full random"
$ws.Range("C9").Value = "sentences"
$ws.Range("D9").Value = 3000
$ws.Range("E9").Value = 200
$ws.Range("F9").Value = 5
$ws.Range("H9").Value = "2000s (ca)"
$ws.Range("I9").Value = "null"
$ws.Range("J9").Value = "yes"
$ws.Range("K9").Value = 115

# Update selection to reflect the saved state
$ws.Range("K14").Select()
